$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update existing cells on Shared_formula sheet ---
$ws1.Range("A3").Value = 3.2323
$ws1.Range("D2").NumberFormat = "0.0000"

# --- New string / bool data block (order controls shared-string table indices) ---
$ws1.Range("B10").Value = "TEXT"
$ws1.Range("B11").Value = "more text"
$ws1.Range("B12").Value = "again more tewt"

$ws1.Range("C10").Value = "ER"
$ws1.Range("D10").Value = "erz"

$ws1.Range("C11").Value = "ezr"
$ws1.Range("D11").Value = "zaeze"

$ws1.Range("C12").Value = "Romain"

$ws1.Range("B15").Value = $true
$ws1.Range("B16").Value = $false

# --- Shared formula block G15:G17 + G18 ---
$ws1.Range("G15:G17").Formula = "=A1"
$ws1.Range("G18").Formula = "=A4"

# --- Add second worksheet "Sheet1" (after the last sheet) with a reference to the shared string "Romain" ---
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet1"
$ws2.Range("A1").Value = "Romain"
$ws2.Range("A3").Select()

# --- Leave the original sheet as the active/selected tab with B17 selected ---
$ws1.Activate()
$ws1.Range("B17").Select()
